# Update the lattice multiplication exercises table: replace the contents
# of every cell (new problem numbers, partial products, and lattice digits)
# while preserving the existing run formatting (sz=32) and <w:br/> layout.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Vertical-tab character used by Word to represent a manual line break
# (<w:br/>) inside Range.Text.
$brk = [char]11

# New content for each of the 15 cells, in row-major (reading) order.
# Each entry is: top "AxB" line, partial-product header line, the dashed
# separator, and the two lattice-row labels.
$cellData = @(
    @("71 x 45", "  4    5", "  ----", "7|    |", "1|    |"),
    @("63 x 65", "  6    5", "  ----", "6|    |", "3|    |"),
    @("24 x 89", "  8    9", "  ----", "2|    |", "4|    |"),
    @("21 x 91", "  9    1", "  ----", "2|    |", "1|    |"),
    @("66 x 23", "  2    3", "  ----", "6|    |", "6|    |"),
    @("17 x 64", "  6    4", "  ----", "1|    |", "7|    |"),
    @("88 x 20", "  2    0", "  ----", "8|    |", "8|    |"),
    @("31 x 85", "  8    5", "  ----", "3|    |", "1|    |"),
    @("56 x 20", "  2    0", "  ----", "5|    |", "6|    |"),
    @("23 x 50", "  5    0", "  ----", "2|    |", "3|    |"),
    @("71 x 85", "  8    5", "  ----", "7|    |", "1|    |"),
    @("17 x 30", "  3    0", "  ----", "1|    |", "7|    |"),
    @("92 x 12", "  1    2", "  ----", "9|    |", "2|    |"),
    @("13 x 56", "  5    6", "  ----", "1|    |", "3|    |"),
    @("52 x 81", "  8    1", "  ----", "5|    |", "2|    |")
)

$numCols = 3
$index = 0
foreach ($lines in $cellData) {
    $row = [math]::Floor($index / $numCols) + 1
    $col = ($index % $numCols) + 1

    $cell = $t.Cell($row, $col)
    $newText = $lines[0] + $brk + $lines[1] + $brk + $lines[2] + $brk + $lines[3] + $brk + $lines[4]
    $cell.Range.Text = $newText

    $index = $index + 1
}
